$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 is currently empty (data starts at row 2), so just populate it directly
# with the new header labels - no row insertion/shifting is needed.

# Fill in the new header row
$ws.Range("A1").Value = "Name of course"
$ws.Range("B1").Value = "Day of Week"
$ws.Range("C1").Value = "Local Time"
$ws.Range("D1").Value = "Number of hours"
$ws.Range("E1").Value = "Start Date"

# Apply the existing bold header style (style index 1 / font id 1) to the header row
$ws.Range("A1:E1").Font.Name = "Calibri"
$ws.Range("A1:E1").Font.Size = 16
$ws.Range("A1:E1").Font.Bold = $true

# Update the "Total:" label to "Total hours:"
$ws.Range("A8").Value = "Total hours:"

# Widen columns to account for the new, wider header text (values chosen so the
# resulting stored column width lands as close as possible to the authored widths)
$ws.Columns.Item(2).ColumnWidth = 17.8125
$ws.Columns.Item(3).ColumnWidth = 15.625
$ws.Columns.Item(4).ColumnWidth = 24.21875
$ws.Columns.Item(5).ColumnWidth = 14.53125
